$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

$para1 = $tr.Paragraphs(1,1)
$oldRun2 = " " + [char]8211 + " pulses high during last cycle of Wait4Instr, Wait4DataWrite or Wait4ReadData, depending upon instruction type. This should move the FSM into the "
$newRun2 = " " + [char]8211 + " pulses high during last cycle of Wait4Instr, Wait4DataWrite or Waitmove4ReadData, depending upon instruction type. This should the FSM into the "
$run2Start = $para1.Start + 20
$sub = $tr.Characters($run2Start, $oldRun2.Length)
Write-Output ("run2 before=[" + $sub.Text + "]")
$sub.Text = $newRun2
Write-Output ("para1 after=[" + $para1.Text + "]")
